# Update UI labels: replace account/review sample data in Sheet1 (A2:B3)
# and append additional scraped review rows (A4:C13), removing the need for
# unused placeholder values. Column C (label) stays blank for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'rahmaawf'
$ws.Range("B2").Value = 'Produk yang saya beli kualitas bagus memuaskan pengemasan cepat, pengiriman cepat, barang tiba sesuai, tidak kurang apa pun, variannya juga pas. Jumlahny sesuai. Harganya murah pokonya mantul'

$ws.Range("A3").Value = 'rafikaameli95'
$ws.Range("B3").Value = 'Performa: barang nya bagus halus dan warna nya cantik banget Kualitas: untuk jahitan rapi banget Cocok Untuk: cocok banget buat di ajak santai santai Jujur yah gua sempat ragu mau pesan boneka di Shopee kek gini tapi ini di luar ekspektasi gua donkk sumpah boneka nya keren bangetttt,warna nya juga cantik banget dan harga nya murah bangettttt,dan akan jadi langganan sih nih toko'

$ws.Range("A4").Value = 'diina_08'
$ws.Range("B4").Value = 'Tekstur: lembut banget, dan nyaman banget dipeluk Desain: imut dan lucu Cocok untuk: anak2 dan orang dewasa tetap cocok Si adek suka banget sama bonekanya, dibawa kemana2 sampe tidur pun sama si pinguin, padahal si adek cowo tapi suka banget sama bonekanya 🥰'
$ws.Range("C4").Style = "Normal"

$ws.Range("A5").Value = 'riskianz_'
$ws.Range("B5").Value = 'Performa: LUCU BANGETTT Kualitas: SEMPURNA Cocok Untuk: BERBAGAI USIA ADUH GATAU LAGI POKOKNYA LUCU BANGET HUAA AMPE NANGIS GATAHAN KARNA GEMOY BANGETT DAN LEMBUT GA KECEWA DEH POKONYA'
$ws.Range("C5").Style = "Normal"

$ws.Range("A6").Value = 'erlindatria30'
$ws.Range("B6").Value = 'Performa: Lucu Cocok Untuk: Semua orang Kualitas: Bagus banget Ini boneka terlucu, terlembut, tergemoy yang aku punya🤩
Produk miniso emang gak usah diraguin lagi
Ini bagus banget, gemoy banget bonekanya😍'
$ws.Range("C6").Style = "Normal"

$ws.Range("A7").Value = 'shereenalicia'
$ws.Range("B7").Value = 'Cocok Untuk: Semua umur Performa: Bagusss Kualitas: Bagusss sesuai harga Boneka sampai dengan aman, packing rapi dan aman, boneka dilapisi plastik bening, penguin nya lucuuu dan soft bangettt, produk miniso ga pernah jelek. Saya udh co produk miniso berkali kali ga pernah kecewaaa'
$ws.Range("C7").Style = "Normal"

$ws.Range("A8").Value = 'ayusetianingrum01'
$ws.Range("B8").Value = 'Performa: lucu Kualitas: bagus Cocok Untuk: kado bagus banget bahannya, lembut dan lucu cocok banget untuk kado 
thankyou seller'
$ws.Range("C8").Style = "Normal"

$ws.Range("A9").Value = 'qoriah19'
$ws.Range("B9").Value = 'Miniso official ya  ngga diragukan lagi,bahannya super lembut,lucu² n affordable.adek suka banget sm pinguinnya..'
$ws.Range("C9").Style = "Normal"

$ws.Range("A10").Value = 'dedijuventini1987'
$ws.Range("B10").Value = 'Kualitas: mantap Cocok Untuk: anak gadis Performa: terbaiklah Alhamdulillah anak gadis qu senang banget,,'
$ws.Range("C10").Style = "Normal"

$ws.Range("A11").Value = 'momikiatar123'
$ws.Range("B11").Value = 'Performa: bagusssssss bgttt😍😍 Kualitas: bagussss Cocok Untuk: kado Mauu nangissss ini baguss bgtt dongg...boneka nya lucuu bgtt..bulunya haluss bgt,ORI minisio..dan udah SNI jg....sukaaaaakkk😍😍😍😍semogaa yg dikasih nya jg sukaa😍😍'
$ws.Range("C11").Style = "Normal"

$ws.Range("A12").Value = 'videvial'
$ws.Range("B12").Value = 'Performa: lembut Cocok Untuk: kado Kualitas: bagus Dapet harga murah karena diskon live, worth banget dengan harga segitu dapet boneka bagus lembut lucu rapih juga mantap deh pokok nya'
$ws.Range("C12").Style = "Normal"

$ws.Range("A13").Value = 'mohamadrichello'
$ws.Range("B13").Value = 'Tekstur: bonekanya lembut Desain: lucu banget Cocok untuk: semua Tx seller. Boneka dah sampai. Lucu banget, bahan juga halus dan lembut. Pas banget dapet harga promo. Respon seller juga cepet, pengiriman juga cepet.'
$ws.Range("C13").Style = "Normal"
